# Apply updated dSF (column F) values after repulling data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    3  = 0
    4  = 0
    6  = 3
    7  = -3
    9  = 0
    12 = 0
    15 = 2
    17 = 0
    19 = 2
    22 = 1
    24 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
